# Auto-generated edit script: updates crypto price/volume table values
# per the commit "Updated cryptos list on Mon Jan 29 11:24:29 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so Excel's COM layer does not
# silently reinterpret numeric-looking strings (e.g. "308.46", "1.00") as
# numbers -- the source data stores these as plain inline text.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.277.78"
$ws.Range("D3").Value = "2.273.32"
$ws.Range("D5").Value = "308.46"
$ws.Range("D6").Value = "97.57"
$ws.Range("D7").Value = "0.526"
$ws.Range("D9").Value = "0.489"
$ws.Range("D10").Value = "34.96"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D13").Value = "6.83"
$ws.Range("D14").Value = "2.625.58"
$ws.Range("D15").Value = "14.60"
$ws.Range("D16").Value = "2.265.45"
$ws.Range("D17").Value = "0.788"
$ws.Range("D18").Value = "42.173.20"
$ws.Range("D19").Value = "12.26"
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("D21").Value = "5.96"
$ws.Range("D22").Value = "67.62"
$ws.Range("D23").Value = "236.67"
$ws.Range("D24").Value = "2.58"
$ws.Range("D26").Value = "1.00"
$ws.Range("D27").Value = "23.59"
$ws.Range("D28").Value = "37.29"
$ws.Range("D29").Value = "9.55"
$ws.Range("D31").Value = "163.53"
$ws.Range("D32").Value = "5.24"
$ws.Range("D35").Value = "17.66"
$ws.Range("D39").Value = "0.115"
$ws.Range("D40").Value = "1.82"
$ws.Range("D41").Value = "4.17"
$ws.Range("D42").Value = "2.27"
$ws.Range("D43").Value = "1.945.33"
$ws.Range("D44").Value = "0.0283"
$ws.Range("D45").Value = "18.71"
$ws.Range("D47").Value = "9.79"
$ws.Range("D48").Value = "54.18"
$ws.Range("D49").Value = "2.498.31"
$ws.Range("D50").Value = "92.06"
$ws.Range("D51").Value = "71.55"

# Restore original (default) cell formatting now that the text values are set.
$priceRange.ClearFormats()

# Coin name / link / volume columns -- plain text updates.
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -1.54%  "
